# "Generate Report for Archive" — refresh the localization-status report:
# the handoff status moves from "Ready for handoff" to "In Translation",
# and the Status column(s) narrow to fit the new (shorter) text.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: per-language status lives in columns E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Replace($oldStatus, $newStatus) | Out-Null
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status lives in column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Replace($oldStatus, $newStatus) | Out-Null
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status lives in column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Replace($oldStatus, $newStatus) | Out-Null
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
